$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 252. This shifts the existing
# row 252 (and everything below it, through row 394) down by one row,
# turning them into rows 253-395, and the overall used range grows
# from A1:R394 to A1:R395.
$ws.Rows("252").Insert()

# Populate the newly inserted row 252 with the new data record.
$ws.Cells.Item(252, 1).Value = 5
$ws.Cells.Item(252, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(252, 3).Value = "Maule"
$ws.Cells.Item(252, 4).Value = 45001
$ws.Cells.Item(252, 5).Value = 7
$ws.Cells.Item(252, 6).Value = 100112009
$ws.Cells.Item(252, 7).Value = "Acelga"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 2500
$ws.Cells.Item(252, 11).Value = 3000
$ws.Cells.Item(252, 12).Value = 3000
$ws.Cells.Item(252, 13).Value = 3000
$ws.Cells.Item(252, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(252, 15).Value = "Región del Maule"
$ws.Cells.Item(252, 16).Value = 750
$ws.Cells.Item(252, 17).Value = 4
$ws.Cells.Item(252, 18).Value = "Hortaliza"
